$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.17200756072998
$ws.Range("B1").Value = 2.890517711639404
$ws.Range("C1").Value = 2.615802049636841
$ws.Range("D1").Value = 2.931085586547852
$ws.Range("E1").Value = 2.878523588180542
